$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.20209999999999
$ws.Range("A8").Value = -22.25820000000002
$ws.Range("A10").Value = -21.6091
$ws.Range("C11").Value = -12.4835
$ws.Range("A12").Value = -21.53299999999999
$ws.Range("C12").Value = -10.2577
$ws.Range("C15").Value = -14.25669999999999
$ws.Range("C17").Value = -13.74989999999999
$ws.Range("A18").Value = -22.06080000000001
$ws.Range("A25").Value = -21.78959999999999
$ws.Range("C26").Value = -12.65640000000001
$ws.Range("C27").Value = -12.9646
$ws.Range("C28").Value = -13.3357
$ws.Range("C32").Value = -12.84160000000001
$ws.Range("A37").Value = -19.56689999999999
$ws.Range("C37").Value = -12.7131
$ws.Range("C41").Value = -12.62660000000001
$ws.Range("C47").Value = -12.4425
$ws.Range("C51").Value = -11.97959999999999
$ws.Range("A55").Value = -22.24939999999999
$ws.Range("C65").Value = -12.4151
$ws.Range("A68").Value = -21.46359999999998
$ws.Range("C73").Value = -10.63630000000001
$ws.Range("A77").Value = -20.37179999999999
$ws.Range("A78").Value = -19.53729999999997
$ws.Range("A79").Value = -20.37079999999998
$ws.Range("A80").Value = -19.54190000000001
$ws.Range("A81").Value = -21.9747
$ws.Range("A82").Value = -21.86710000000002
$ws.Range("A84").Value = -21.36749999999998
$ws.Range("C84").Value = -12.8098
$ws.Range("C85").Value = -12.8786
$ws.Range("C89").Value = -13.8461
$ws.Range("C93").Value = -10.34119999999999
$ws.Range("C95").Value = -12.4258
$ws.Range("C98").Value = -12.77210000000002
$ws.Range("C99").Value = -12.1021
$ws.Range("A101").Value = -20.23809999999999
$ws.Range("C101").Value = -12.2978
$ws.Range("A102").Value = -20.40319999999998
$ws.Range("C102").Value = -12.3923
